$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductDetails")
$ws.Activate()

# Update B2 value from 1 to 2
$ws.Range("B2").Value = 2

# Update C2 shared string from "XS" to "L"
$ws.Range("C2").Value = "L"

# Update selection on the sheet
$ws.Range("N6").Select()
